$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H43").Value = 4569.1
$ws.Range("I43").Value = 4340.2
$ws.Range("J43").Value = 4798
$ws.Range("K43").Value = 4340.2
$ws.Range("L43").Value = 4798
$ws.Range("M43").Value = -4271.2
$ws.Range("N43").Value = -4936

$ws.Range("H53").Value = 976.0769
$ws.Range("I53").Value = 753.1
$ws.Range("J53").Value = 1719.3334
$ws.Range("K53").Value = 753.1
$ws.Range("L53").Value = 1719.3334
$ws.Range("M53").Value = -116.1
$ws.Range("N53").Value = -2993.3334

$ws.Range("H64").Value = 7291.4165
$ws.Range("I64").Value = 4833.3335
$ws.Range("J64").Value = 9749.5
$ws.Range("K64").Value = 4833.3335
$ws.Range("L64").Value = 9749.5
$ws.Range("M64").Value = -4585.3335
$ws.Range("N64").Value = -10245.5

$ws.Range("H67").Value = 7291.4165
$ws.Range("I67").Value = 4833.3335
$ws.Range("J67").Value = 9749.5
$ws.Range("K67").Value = 4833.3335
$ws.Range("L67").Value = 9749.5
$ws.Range("M67").Value = -3975.3335
$ws.Range("N67").Value = -11465.5

$ws.Range("H137").Value = 836.7778
$ws.Range("I137").Value = 788.625
$ws.Range("J137").Value = 1222
$ws.Range("K137").Value = 2365.875
$ws.Range("L137").Value = 3666
$ws.Range("M137").Value = 184.125
$ws.Range("N137").Value = -8766

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H5").Value = 2688.3333
$ws.Range("I5").Value = 1695.3334
$ws.Range("J5").Value = 3681.3333
$ws.Range("K5").Value = 1695.3334
$ws.Range("L5").Value = 3681.3333
$ws.Range("M5").Value = -1583.3334
$ws.Range("N5").Value = -3905.3333

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H121").Value = 99999
$ws.Range("J121").Value = 99999
$ws.Range("L121").Value = 99999
$ws.Range("N121").Value = -103493

$ws.Range("H132").Value = 601.25
$ws.Range("I132").Value = 636.4545000000001
$ws.Range("J132").Value = 214
$ws.Range("K132").Value = 1909.3635
$ws.Range("L132").Value = 642
$ws.Range("M132").Value = 620.6364999999998
$ws.Range("N132").Value = -5702

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 2688.3333
$ws.Range("I4").Value = 1695.3334
$ws.Range("J4").Value = 3681.3333
$ws.Range("K4").Value = 1695.3334
$ws.Range("L4").Value = 3681.3333
$ws.Range("M4").Value = -1580.3334
$ws.Range("N4").Value = -3911.3333

$ws.Range("H132").Value = 80000
$ws.Range("J132").Value = 80000
$ws.Range("L132").Value = 80000
$ws.Range("N132").Value = -90120

$ws.Range("H134").Value = 1267.4286
$ws.Range("I134").Value = 1334.1538
$ws.Range("K134").Value = 4002.4614
$ws.Range("M134").Value = -1467.4614

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H132").Value = 1995.591
$ws.Range("I132").Value = 2035.619
$ws.Range("J132").Value = 1155
$ws.Range("K132").Value = 6106.857
$ws.Range("L132").Value = 3465
$ws.Range("M132").Value = -3576.857
$ws.Range("N132").Value = -8525

$ws.Range("H134").Value = 1449.75
$ws.Range("I134").Value = 1449.75
$ws.Range("K134").Value = 4349.25
$ws.Range("M134").Value = -1814.25

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 1500
$ws.Range("M3").Value = -1388

$ws.Range("H114").Value = 420
$ws.Range("I114").Value = 420
$ws.Range("K114").Value = 1260
$ws.Range("M114").Value = 1994

$ws.Range("H122").Value = 703.2857
$ws.Range("I122").Value = 566.6667
$ws.Range("J122").Value = 805.75
$ws.Range("K122").Value = 5100.0003
$ws.Range("L122").Value = 7251.75
$ws.Range("M122").Value = -2650.0003
$ws.Range("N122").Value = -12151.75

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H102").Value = 2352.4443
$ws.Range("I102").Value = 2209
$ws.Range("K102").Value = 2209
$ws.Range("M102").Value = -587

$ws.Range("H132").Value = 4649.75
$ws.Range("I132").Value = 4649.75
$ws.Range("K132").Value = 13949.25
$ws.Range("M132").Value = -11419.25

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 3497
$ws.Range("I7").Value = 1994
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 1994
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -1882
$ws.Range("N7").Value = -5224

$ws.Range("H61").Value = 969.125
$ws.Range("I61").Value = 893.2857
$ws.Range("K61").Value = 893.2857
$ws.Range("M61").Value = -691.2857

$ws.Range("H68").Value = 2439.6
$ws.Range("I68").Value = 2550
$ws.Range("J68").Value = 1998
$ws.Range("K68").Value = 2550
$ws.Range("L68").Value = 1998
$ws.Range("M68").Value = -1801
$ws.Range("N68").Value = -3496

$ws.Range("H71").Value = 2439.6
$ws.Range("I71").Value = 2550
$ws.Range("J71").Value = 1998
$ws.Range("K71").Value = 12750
$ws.Range("L71").Value = 9990
$ws.Range("M71").Value = -9006
$ws.Range("N71").Value = -17478

$ws.Range("H113").Value = 969.125
$ws.Range("I113").Value = 893.2857
$ws.Range("K113").Value = 893.2857
$ws.Range("M113").Value = 1276.7143

$ws.Range("H126").Value = 3497
$ws.Range("I126").Value = 1994
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5982
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3512
$ws.Range("N126").Value = -19940

$ws.Range("H136").Value = 1712.25
$ws.Range("I136").Value = 1712.25
$ws.Range("K136").Value = 5136.75
$ws.Range("M136").Value = -2586.75

$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 7214.143
$ws.Range("I62").Value = 7099.4
$ws.Range("J62").Value = 7501
$ws.Range("K62").Value = 7099.4
$ws.Range("L62").Value = 7501
$ws.Range("M62").Value = -6475.4
$ws.Range("N62").Value = -8749

$ws.Range("H65").Value = 7214.143
$ws.Range("I65").Value = 7099.4
$ws.Range("J65").Value = 7501
$ws.Range("K65").Value = 35497
$ws.Range("L65").Value = 37505
$ws.Range("M65").Value = -32377
$ws.Range("N65").Value = -43745
